$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the Neo4j query text in cell B4 (shared string index 13)
$lines = @(
    "MATCH (f:file)-->(s:study)",
    "OPTIONAL MATCH (samp:sample)<--(f)",
    "OPTIONAL MATCH (samp)-->(p:participant)",
    "OPTIONAL MATCH (f)<--(g:genomic_info)",
    "OPTIONAL MATCH (p)<--(diag:diagnosis)",
    "WITH s, p, samp, f, g, diag",
    "WHERE  f.file_type in ['BAM']",
    "WITH DISTINCT f, s, p, samp",
    "RETURN",
    "    coalesce(f.file_name, '') as ``File Name``,",
    "    coalesce(s.study_name,'') as ``Study Name``,",
    "    coalesce(s.phs_accession,'') as ``Accession``,",
    "    coalesce(p.participant_id, '') as ``Participant ID``,",
    "    coalesce(samp.sample_id, '') as ``Sample ID``,",
    "    coalesce(f.file_type, '') as ``File Type``",
    "ORDER BY f.file_name limit 100"
)
$newQuery = [string]::Join([char]10, $lines)

$ws.Range("B4").Value = $newQuery

# Update row 4 height to match new (taller) wrapped text
$ws.Rows.Item(4).RowHeight = 248

# Update the selected cell on the sheet
$ws.Range("B5").Select()
